# Weekly update: insert two new price records (rows 9 and 10) for
# "Arveja Verde" at Terminal Hortofrutícola Agro Chillán, shifting the
# existing historical rows down by two positions.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows above the current row 9 so the previous rows
# 9..58 become rows 11..60.
$ws.Rows.Item(9).Insert()
$ws.Rows.Item(9).Insert()

# New row 9
$ws.Cells.Item(9, 1).Value = 7
$ws.Cells.Item(9, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(9, 3).Value = "Ñuble"
$ws.Cells.Item(9, 4).Value = 44592
$ws.Cells.Item(9, 5).Value = 16
$ws.Cells.Item(9, 6).Value = 100112022
$ws.Cells.Item(9, 7).Value = "Arveja Verde"
$ws.Cells.Item(9, 8).Value = "Sin especificar"
$ws.Cells.Item(9, 9).Value = "Primera"
$ws.Cells.Item(9, 10).Value = 30
$ws.Cells.Item(9, 11).Value = 24000
$ws.Cells.Item(9, 12).Value = 24000
$ws.Cells.Item(9, 13).Value = 24000
$ws.Cells.Item(9, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(9, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(9, 16).Value = 960
$ws.Cells.Item(9, 17).Value = 25
$ws.Cells.Item(9, 18).Value = "Hortaliza"

# New row 10
$ws.Cells.Item(10, 1).Value = 7
$ws.Cells.Item(10, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(10, 3).Value = "Ñuble"
$ws.Cells.Item(10, 4).Value = 44592
$ws.Cells.Item(10, 5).Value = 16
$ws.Cells.Item(10, 6).Value = 100112022
$ws.Cells.Item(10, 7).Value = "Arveja Verde"
$ws.Cells.Item(10, 8).Value = "Sin especificar"
$ws.Cells.Item(10, 9).Value = "Segunda"
$ws.Cells.Item(10, 10).Value = 30
$ws.Cells.Item(10, 11).Value = 23000
$ws.Cells.Item(10, 12).Value = 23000
$ws.Cells.Item(10, 13).Value = 23000
$ws.Cells.Item(10, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(10, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(10, 16).Value = 920
$ws.Cells.Item(10, 17).Value = 25
$ws.Cells.Item(10, 18).Value = "Hortaliza"
